$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1659919028340081
$ws.Range("C2").Value = 0.5991902834008097
$ws.Range("J2").Value = 0.008097165991902834
$ws.Range("P2").Value = 0.1295546558704453
$ws.Range("S2").Value = 0.09716599190283401
$ws.Range("C3").Value = 0.01351351351351351
$ws.Range("J3").Value = 0.006756756756756757
$ws.Range("P3").Value = 0.7905405405405406
$ws.Range("S3").Value = 0.1891891891891892
$ws.Range("J4").Value = 0.025
$ws.Range("P4").Value = 0.725
$ws.Range("S4").Value = 0.25
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.03097345132743363
$ws.Range("D6").Value = 0.01769911504424779
$ws.Range("F6").Value = 0.05309734513274336
$ws.Range("J6").Value = 0.2079646017699115
$ws.Range("O6").Value = 0.01769911504424779
$ws.Range("Q6").Value = 0.1725663716814159
$ws.Range("R6").Value = 0.084070796460177
$ws.Range("S6").Value = 0.415929203539823
$ws.Range("B7").Value = 0.09836065573770492
$ws.Range("D7").Value = 0.03278688524590164
$ws.Range("F7").Value = 0.06010928961748634
$ws.Range("J7").Value = 0.1256830601092896
$ws.Range("O7").Value = 0.01092896174863388
$ws.Range("Q7").Value = 0.1748633879781421
$ws.Range("R7").Value = 0.07650273224043716
$ws.Range("S7").Value = 0.4207650273224044
$ws.Range("B8").Value = 0.07692307692307693
$ws.Range("D8").Value = 0.01923076923076923
$ws.Range("F8").Value = 0.05555555555555555
$ws.Range("J8").Value = 0.1239316239316239
$ws.Range("O8").Value = 0.01923076923076923
$ws.Range("Q8").Value = 0.2136752136752137
$ws.Range("R8").Value = 0.09188034188034189
$ws.Range("S8").Value = 0.3995726495726496
$ws.Range("B9").Value = 0.09216589861751152
$ws.Range("D9").Value = 0.004608294930875576
$ws.Range("F9").Value = 0.1059907834101382
$ws.Range("J9").Value = 0.08755760368663594
$ws.Range("O9").Value = 0.02304147465437788
$ws.Range("Q9").Value = 0.1797235023041475
$ws.Range("R9").Value = 0.08294930875576037
$ws.Range("S9").Value = 0.423963133640553
$ws.Range("B10").Value = 0.1005770816158285
$ws.Range("D10").Value = 0.01731244847485573
$ws.Range("E10").Value = 0.00247320692497939
$ws.Range("F10").Value = 0.0717230008244023
$ws.Range("J10").Value = 0.1005770816158285
$ws.Range("O10").Value = 0.00989282769991756
$ws.Range("Q10").Value = 0.2324814509480627
$ws.Range("R10").Value = 0.1005770816158285
$ws.Range("S10").Value = 0.3643858202802968
$ws.Range("G11").Value = 0.1574803149606299
$ws.Range("J11").Value = 0.05511811023622047
$ws.Range("K11").Value = 0.1811023622047244
$ws.Range("L11").Value = 0.594488188976378
$ws.Range("S11").Value = 0.01181102362204724
$ws.Range("G12").Value = 0.7468354430379747
$ws.Range("J12").Value = 0.1708860759493671
$ws.Range("K12").Value = 0.01265822784810127
$ws.Range("L12").Value = 0.04430379746835443
$ws.Range("S12").Value = 0.02531645569620253
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("S13").Value = 0.1111111111111111
$ws.Range("F15").Value = 0.01895734597156398
$ws.Range("H15").Value = 0.1753554502369668
$ws.Range("I15").Value = 0.06161137440758294
$ws.Range("J15").Value = 0.3554502369668247
$ws.Range("K15").Value = 0.06161137440758294
$ws.Range("M15").Value = 0.01895734597156398
$ws.Range("O15").Value = 0.07582938388625593
$ws.Range("S15").Value = 0.2322274881516588
$ws.Range("F16").Value = 0.0335195530726257
$ws.Range("H16").Value = 0.1955307262569832
$ws.Range("I16").Value = 0.09497206703910614
$ws.Range("J16").Value = 0.3798882681564246
$ws.Range("K16").Value = 0.1340782122905028
$ws.Range("M16").Value = 0.0223463687150838
$ws.Range("O16").Value = 0.02793296089385475
$ws.Range("S16").Value = 0.111731843575419
$ws.Range("F17").Value = 0.006036217303822937
$ws.Range("H17").Value = 0.1549295774647887
$ws.Range("I17").Value = 0.096579476861167
$ws.Range("J17").Value = 0.4507042253521127
$ws.Range("K17").Value = 0.08249496981891348
$ws.Range("M17").Value = 0.02414486921529175
$ws.Range("O17").Value = 0.0744466800804829
$ws.Range("S17").Value = 0.1106639839034205
$ws.Range("F18").Value = 0.02714932126696833
$ws.Range("H18").Value = 0.1809954751131222
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.416289592760181
$ws.Range("K18").Value = 0.09502262443438914
$ws.Range("M18").Value = 0.01357466063348416
$ws.Range("O18").Value = 0.09049773755656108
$ws.Range("S18").Value = 0.09954751131221719
$ws.Range("F19").Value = 0.01462225832656377
$ws.Range("H19").Value = 0.2250203086921202
$ws.Range("I19").Value = 0.09991876523151909
$ws.Range("J19").Value = 0.3574329813160033
$ws.Range("K19").Value = 0.08610885458976442
$ws.Range("M19").Value = 0.02437043054427295
$ws.Range("N19").Value = 0.0008123476848090983
$ws.Range("O19").Value = 0.06823720552396426
$ws.Range("S19").Value = 0.1234768480909829
